# Updated cryptos list on Sat Jun 10 09:47:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking strings (e.g. "238.49", "4.440",
# "0.000006357") as text (inlineStr) in the source data, so force the whole
# column to a text format before writing; otherwise Excel's COM layer would
# auto-coerce the assigned strings into real numbers and silently normalize
# (e.g. "4.440" -> 4.44, "1.000" -> 1, "0.000006357" -> 6.357E-06).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.687.31"
$ws.Range("E2").Value = "  -3.58%  "
$ws.Range("D3").Value = "1.739.34"
$ws.Range("E3").Value = "  -5.72%  "
$ws.Range("D5").Value = "238.49"
$ws.Range("E5").Value = "  -8.32%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.4906"
$ws.Range("E7").Value = "  -6.86%  "
$ws.Range("E8").Value = "  -6.93%  "
$ws.Range("D9").Value = "0.2404"
$ws.Range("E9").Value = "  -23.87%  "
$ws.Range("D10").Value = "0.05993"
$ws.Range("E10").Value = "  -11.84%  "
$ws.Range("D11").Value = "1.747.78"
$ws.Range("E11").Value = "  -5.16%  "
$ws.Range("D12").Value = "0.06782"
$ws.Range("E12").Value = "  -12.68%  "
$ws.Range("D13").Value = "14.69"
$ws.Range("E13").Value = "  -21.76%  "
$ws.Range("D14").Value = "4.440"
$ws.Range("E14").Value = "  -11.51%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "76.41"
$ws.Range("E15").Value = "  -13.26%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.5806"
$ws.Range("E16").Value = "  -25.97%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "25.716.77"
$ws.Range("E20").Value = "  -17.47%  "
$ws.Range("D21").Value = "0.000006357"
$ws.Range("E21").Value = "  -19.82%  "
$ws.Range("D22").Value = "1.957.69"
$ws.Range("E22").Value = "  -5.74%  "
$ws.Range("D23").Value = "3.932"
$ws.Range("E23").Value = "  -14.73%  "
$ws.Range("D24").Value = "5.087"
$ws.Range("E24").Value = "  -14.97%  "
$ws.Range("D25").Value = "7.826"
$ws.Range("E25").Value = "  -16.02%  "
$ws.Range("D26").Value = "136.16"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").Value = "1.458"
$ws.Range("E27").Value = "  -13.13%  "
$ws.Range("D28").Value = "1.835"
$ws.Range("E28").Value = "  -17.35%  "
$ws.Range("D29").Value = "14.44"
$ws.Range("E29").Value = "  -15.24%  "
$ws.Range("D30").Value = "99.86"
$ws.Range("E30").Value = "  -10.08%  "
$ws.Range("D31").Value = "0.08094"
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("D32").Value = "3.716"
$ws.Range("E32").Value = "  -11.46%  "
$ws.Range("D33").Value = "3.349"
$ws.Range("E33").Value = "  -17.96%  "
$ws.Range("D34").Value = "0.04360"
$ws.Range("E34").Value = "  -10.74%  "
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "2.687"
$ws.Range("E36").Value = "  -5.92%  "
$ws.Range("D37").Value = "1.018"
$ws.Range("E37").Value = "  -10.86%  "
$ws.Range("D38").Value = "0.5975"
$ws.Range("E38").Value = "  -18.14%  "
$ws.Range("D39").Value = "2.726"
$ws.Range("E39").Value = "  -11.98%  "
$ws.Range("D40").Value = "2.039"
$ws.Range("E40").Value = "  -10.54%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "102.73"
$ws.Range("E42").Value = "  -6.52%  "
$ws.Range("E43").Value = "  -14.45%  "
$ws.Range("D44").Value = "0.7831"
$ws.Range("E44").Value = "  -12.99%  "
$ws.Range("D45").Value = "0.3786"
$ws.Range("E45").Value = "  -20.77%  "
$ws.Range("D46").Value = "5.120"
$ws.Range("E46").Value = "  -13.97%  "
$ws.Range("D47").Value = "5.989"
$ws.Range("E47").Value = "  -22.15%  "
$ws.Range("E48").Value = "  -12.45%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1062"
$ws.Range("E49").Value = "  -14.21%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "30.10"
$ws.Range("E50").Value = "  -13.45%  "
$ws.Range("D51").Value = "52.13"
$ws.Range("E51").Value = "  -12.90%  "
